$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----------------------------------------------------------------------
# 1) Cell value updates (table renamed, columns reordered/renamed,
#    attributes updated, index rows populated)
# ----------------------------------------------------------------------

$ws.Range("B1").Value2  = "Pagamentos"

$ws.Range("A5").Value2  = "id"
$ws.Range("A6").Value2  = "movimentacao_id"

$ws.Range("A7").Value2  = "metodo_pagamento"
$ws.Range("C7").Value2  = "ENUM('Cartão', 'Dinheiro', 'Pix')"
$ws.Range("D7").Value2  = "3 – 8"
$ws.Range("E7").Value2  = "NOT NULL"
$ws.Range("H7").Value2  = "Forma de pagamento utilizada"

$ws.Range("A8").Value2  = "valor_pago"
$ws.Range("C8").Value2  = "decimal(10, 2)"
$ws.Range("D8").Value2  = "0 - sem limite"
$ws.Range("E8").Value2  = "NOT NULL"
$ws.Range("H8").Value2  = "Valor do pagamento cobrado pelo estacionamento"

$ws.Range("A9").Value2  = "data_pagamento"
$ws.Range("C9").Value2  = "timestamp"
$ws.Range("D9").Value2  = "sem limite"
$ws.Range("E9").Value2  = "NOT NULL"
$ws.Range("H9").Value2  = "Data e hora que foi realizado pagamento"

$ws.Range("A13").Value2 = "PRIMARY"
$ws.Range("C13").Value2 = "Sim"
$ws.Range("D13").Value2 = "Não"
$ws.Range("E13").Value2 = "Sim"
$ws.Range("F13").Value2 = "id"

$ws.Range("A14").Value2 = "INDEX_movimentacao_id"
$ws.Range("C14").Value2 = "Não"
$ws.Range("D14").Value2 = "Sim"
$ws.Range("E14").Value2 = "Não"
$ws.Range("F14").Value2 = "movimentacao_id"

# ----------------------------------------------------------------------
# 2) Re-apply formatting so every cell carries the (reordered) style that
#    the target workbook expects. Every needed look already exists
#    somewhere in the sheet, so copy it across with PasteSpecial(formats).
#    NOTE: D9's own look changes (picks up D5's), so grab D9's original
#    look for D7 *before* touching D9.
# ----------------------------------------------------------------------

$xlPasteFormats = -4122

function Copy-Format($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($dstAddr).PasteSpecial($xlPasteFormats) | Out-Null
}

# D7 must borrow D9's pristine (pre-edit) look first.
Copy-Format "D9" "D7"

Copy-Format "A1" "A1:A2,C4:H4,C12:E12"
Copy-Format "B1" "B1:H1,B2:H2"
Copy-Format "B3" "A3:H3"
Copy-Format "B4" "A4:B4,A11:H11,A12:B12,F12:H12"
Copy-Format "A5" "A5:B10"
Copy-Format "D5" "C5:H6,C7,E7:H7,C8:H10"
Copy-Format "A13" "A13,F13,A14,F14,A15:A16,F16"
Copy-Format "B13" "B13,H13,B14,H14,B15:B16,H16"
Copy-Format "C13" "C13:C16"
Copy-Format "D13" "D13:E14,D15:E16"
Copy-Format "G13" "G13:G14,G16"

$excel.CutCopyMode = $false

# ----------------------------------------------------------------------
# 3) Column C gets an explicit (best-fit-like) width, matching the new
#    "Tipo do Dado" column that now holds longer text (ENUM(...) etc).
# ----------------------------------------------------------------------

$ws.Columns("C").ColumnWidth = 26.25

# ----------------------------------------------------------------------
# 4) Selection moves to F15:H15 (last edited cell block).
# ----------------------------------------------------------------------

$ws.Range("F15:H15").Select() | Out-Null
